$p = $ppt.ActivePresentation

# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" / "Office" colour scheme (blue accents)
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet" colour scheme (pink/purple accents),
#                            which is the theme actually wired to the (only) slide master.
# The authored edit swaps the two themes' content, so the presentation's visible design
# becomes the plain "Office" blue palette. The PowerPoint object model only exposes the
# theme that is bound to the slide master/slides (theme2.xml) for read/write, so we
# rewrite its 12 colour-scheme slots, in order, to the Office palette values that
# originally lived in theme1.xml.

$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}

Write-Output "theme colours updated"
